# Update my log.
$wb = $excel.ActiveWorkbook

# Add the new worksheet ("Sheet2") after the existing Sheet1.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Row 1
$ws2.Range("A1").Value = "今天7:31 2011-04-01查我的学车资料状态已经是已结业了"

# Row 4
$ws2.Range("A4").Value = "[查看]"
$ws2.Range("B4").Value = " 张昀"
$ws2.Range("C4").Value = " C1"
$ws2.Range("D4").Value = " 已结业"
$ws2.Range("E4").Value = " 初次申请"
$ws2.Range("F4").Value = " 2010-10-25 00:00:00"
$ws2.Range("J4").Value = " 广州市穗通驾驶员培训有限公司"
$ws2.Range("K4").Value = " JP440106000014"
$ws2.Range("L4").Value = " S10111016"

# Row 5
$ws2.Range("A5").Value = "当前处于阶段"

# Row 6
$ws2.Range("A6").Value = "大纲要求培训时长为:"

# Row 7
$ws2.Range("A7").Value = "理论"

# Row 8
$ws2.Range("A8").Value = "实操"

# Row 9
$ws2.Range("A9").Value = "实际完成："

# Row 10
$ws2.Range("A10").Value = "理论"

# Row 11
$ws2.Range("A11").Value = "实操"

# Row 12
$ws2.Range("A12").Value = "学员名称"
$ws2.Range("B12").Value = "培训车型"
$ws2.Range("C12").Value = "培训状态"
$ws2.Range("D12").Value = "培训类型"
$ws2.Range("E12").Value = "入学时间"

# Row 13
$ws2.Range("A13").Value = 5
$ws2.Range("B13").Value = 0
$ws2.Range("C13").Value = 0
$ws2.Range("D13").Value = 0
$ws2.Range("E13").Value = 0

# Selection / active-sheet bookkeeping to match the target state.
$ws2.Range("A1:L13").Select()
